# Server.xlsx "fix error in linux" edit
# - Replace the placeholder 127.0.0.1 IPs in column F (rows 2-6) with real
#   per-server LAN addresses.
# - Remove the stray duplicate "GameServer_2" row (row 7) that was left
#   over from testing - clear its contents entirely.
# - Column F (IP) is now wide enough to show the longer dotted addresses.
# - Selection cursor ends up parked on the (now blank) row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP column for the five real server rows.
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Wipe out the obsolete extra row (was GameServer_2 / duplicate data).
$ws.Range("A7:H7").ClearContents()

# Widen column F to fit the new IP strings (matches the saved workbook's
# stored column width of 15 characters).
$ws.Columns("F").ColumnWidth = 14.285714285714286

# Select the now-empty row 7, matching the cursor position left behind
# after clearing it.
[void]$ws.Rows(7).Select()
